$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 1.780241
$ws.Range("N2").Value = 5.340723000000001
$ws.Range("O2").Value = 0.2571704734300857
$ws.Range("P2").Value = 0.2571704734300857
$ws.Range("Q2").Value = 118.344999065914
$ws.Range("R2").Value = 1065.104991593226
$ws.Range("S2").Value = 0.01108914736834913
$ws.Range("T2").Value = 0.01108914736834913
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("O3").Value = 0.4607709215973151
$ws.Range("P3").Value = 0.4607709215973152
$ws.Range("Q3").Value = 212.038083372192
$ws.Range("R3").Value = 1908.342750349728
$ws.Range("S3").Value = 0.01986836429739573
$ws.Range("T3").Value = 0.01986836429739573
$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4746316666666666
$ws.Range("N4").Value = 1.423895
$ws.Range("O4").Value = 0.06856445302718973
$ws.Range("P4").Value = 0.06856445302718973
$ws.Range("Q4").Value = 31.55206747194333
$ws.Range("R4").Value = 283.96860724749
$ws.Range("S4").Value = 0.002956487631366668
$ws.Range("T4").Value = 0.002956487631366667
$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("M5").Value = 0.9303213333333334
$ws.Range("N5").Value = 2.790964
$ws.Range("O5").Value = 0.1343925781596098
$ws.Range("P5").Value = 0.1343925781596098
$ws.Range("Q5").Value = 61.84492848121867
$ws.Range("R5").Value = 556.604356330968
$ws.Range("S5").Value = 0.005794985266181595
$ws.Range("T5").Value = 0.005794985266181594
$ws.Range("G6").Value = 66.47695399999999
$ws.Range("H6").Value = 199.430862
$ws.Range("I6").Value = 0.04311983106164722
$ws.Range("J6").Value = 0.04311983106164721
$ws.Range("M6").Value = 0.547574
$ws.Range("N6").Value = 1.642722
$ws.Range("O6").Value = 0.07910157378579964
$ws.Range("P6").Value = 0.07910157378579964
$ws.Range("Q6").Value = 36.401051609596
$ws.Range("R6").Value = 327.609464486364
$ws.Range("S6").Value = 0.003410846498354103
$ws.Range("T6").Value = 0.003410846498354102
$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("M7").Value = 1.780241
$ws.Range("N7").Value = 5.340723000000001
$ws.Range("O7").Value = 0.2571704734300857
$ws.Range("P7").Value = 0.2571704734300857
$ws.Range("Q7").Value = 2423.582835175629
$ws.Range("R7").Value = 21812.24551658066
$ws.Range("S7").Value = 0.227094236603063
$ws.Range("T7").Value = 0.227094236603063
$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("O8").Value = 0.4607709215973151
$ws.Range("P8").Value = 0.4607709215973152
$ws.Range("S8").Value = 0.406883493635124
$ws.Range("T8").Value = 0.406883493635124
$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4746316666666666
$ws.Range("N9").Value = 1.423895
$ws.Range("O9").Value = 0.06856445302718973
$ws.Range("P9").Value = 0.06856445302718973
$ws.Range("Q9").Value = 646.1536164845849
$ws.Range("R9").Value = 5815.382548361265
$ws.Range("S9").Value = 0.06054580026485523
$ws.Range("T9").Value = 0.06054580026485522
$ws.Range("I10").Value = 0.8830494168872806
$ws.Range("J10").Value = 0.8830494168872804
$ws.Range("M10").Value = 0.9303213333333334
$ws.Range("N10").Value = 2.790964
$ws.Range("O10").Value = 0.1343925781596098
$ws.Range("P10").Value = 0.1343925781596098
$ws.Range("Q10").Value = 1266.519990644172
$ws.Range("R10").Value = 11398.67991579755
$ws.Range("S10").Value = 0.1186752877778217
$ws.Range("T10").Value = 0.1186752877778217
$ws.Range("I11").Value = 0.8830494168872806
$ws.Range("J11").Value = 0.8830494168872804
$ws.Range("M11").Value = 0.547574
$ws.Range("N11").Value = 1.642722
$ws.Range("O11").Value = 0.07910157378579964
$ws.Range("P11").Value = 0.07910157378579964
$ws.Range("Q11").Value = 745.4557823286059
$ws.Range("R11").Value = 6709.102040957454
$ws.Range("S11").Value = 0.06985059860641657
$ws.Range("T11").Value = 0.06985059860641657
$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("M12").Value = 1.780241
$ws.Range("N12").Value = 5.340723000000001
$ws.Range("O12").Value = 0.2571704734300857
$ws.Range("P12").Value = 0.2571704734300857
$ws.Range("Q12").Value = 79.810183657992
$ws.Range("R12").Value = 718.2916529219281
$ws.Range("S12").Value = 0.007478363218251022
$ws.Range("T12").Value = 0.007478363218251021
$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("O13").Value = 0.4607709215973151
$ws.Range("P13").Value = 0.4607709215973152
$ws.Range("Q13").Value = 142.995466728576
$ws.Range("R13").Value = 1286.959200557184
$ws.Range("S13").Value = 0.0133989422119634
$ws.Range("T13").Value = 0.0133989422119634
$ws.Range("G14").Value = 44.831112
$ws.Range("H14").Value = 134.493336
$ws.Range("I14").Value = 0.02907940059566787
$ws.Range("J14").Value = 0.02907940059566786
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.4746316666666666
$ws.Range("N14").Value = 1.423895
$ws.Range("O14").Value = 0.06856445302718973
$ws.Range("P14").Value = 0.06856445302718973
$ws.Range("Q14").Value = 21.27826540708
$ws.Range("R14").Value = 191.50438866372
$ws.Range("S14").Value = 0.001993813196200503
$ws.Range("T14").Value = 0.001993813196200503
$ws.Range("G15").Value = 44.831112
$ws.Range("H15").Value = 134.493336
$ws.Range("I15").Value = 0.02907940059566787
$ws.Range("J15").Value = 0.02907940059566786
$ws.Range("M15").Value = 0.9303213333333334
$ws.Range("N15").Value = 2.790964
$ws.Range("O15").Value = 0.1343925781596098
$ws.Range("P15").Value = 0.1343925781596098
$ws.Range("Q15").Value = 41.707339890656
$ws.Range("R15").Value = 375.366059015904
$ws.Range("S15").Value = 0.003908055617387898
$ws.Range("T15").Value = 0.003908055617387897
$ws.Range("G16").Value = 44.831112
$ws.Range("H16").Value = 134.493336
$ws.Range("I16").Value = 0.02907940059566787
$ws.Range("J16").Value = 0.02907940059566786
$ws.Range("M16").Value = 0.547574
$ws.Range("N16").Value = 1.642722
$ws.Range("O16").Value = 0.07910157378579964
$ws.Range("P16").Value = 0.07910157378579964
$ws.Range("Q16").Value = 24.548351322288
$ws.Range("R16").Value = 220.935161900592
$ws.Range("S16").Value = 0.002300226351865048
$ws.Range("T16").Value = 0.002300226351865048
$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("M17").Value = 1.780241
$ws.Range("N17").Value = 5.340723000000001
$ws.Range("O17").Value = 0.2571704734300857
$ws.Range("P17").Value = 0.2571704734300857
$ws.Range("Q17").Value = 94.05123755966102
$ws.Range("R17").Value = 846.4611380369493
$ws.Range("S17").Value = 0.008812776557578141
$ws.Range("T17").Value = 0.00881277655757814
$ws.Range("G18").Value = 52.83062100000001
$ws.Range("H18").Value = 158.491863
$ws.Range("I18").Value = 0.0342682285413064
$ws.Range("J18").Value = 0.03426822854130639
$ws.Range("O18").Value = 0.4607709215973151
$ws.Range("P18").Value = 0.4607709215973152
$ws.Range("Q18").Value = 168.511084611408
$ws.Range("R18").Value = 1516.599761502672
$ws.Range("S18").Value = 0.01578980324648517
$ws.Range("T18").Value = 0.01578980324648517
$ws.Range("G19").Value = 52.83062100000001
$ws.Range("H19").Value = 158.491863
$ws.Range("I19").Value = 0.0342682285413064
$ws.Range("J19").Value = 0.03426822854130639
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.4746316666666666
$ws.Range("N19").Value = 1.423895
$ws.Range("O19").Value = 0.06856445302718973
$ws.Range("P19").Value = 0.06856445302718973
$ws.Range("Q19").Value = 25.075085696265
$ws.Range("R19").Value = 225.675771266385
$ws.Range("S19").Value = 0.002349582346145405
$ws.Range("T19").Value = 0.002349582346145405
$ws.Range("G20").Value = 52.83062100000001
$ws.Range("H20").Value = 158.491863
$ws.Range("I20").Value = 0.0342682285413064
$ws.Range("J20").Value = 0.03426822854130639
$ws.Range("M20").Value = 0.9303213333333334
$ws.Range("N20").Value = 2.790964
$ws.Range("O20").Value = 0.1343925781596098
$ws.Range("P20").Value = 0.1343925781596098
$ws.Range("Q20").Value = 49.14945376954801
$ws.Range("R20").Value = 442.3450839259321
$ws.Range("S20").Value = 0.004605395582628892
$ws.Range("T20").Value = 0.004605395582628891
$ws.Range("G21").Value = 52.83062100000001
$ws.Range("H21").Value = 158.491863
$ws.Range("I21").Value = 0.0342682285413064
$ws.Range("J21").Value = 0.03426822854130639
$ws.Range("M21").Value = 0.547574
$ws.Range("N21").Value = 1.642722
$ws.Range("O21").Value = 0.07910157378579964
$ws.Range("P21").Value = 0.07910157378579964
$ws.Range("Q21").Value = 28.928674463454
$ws.Range("R21").Value = 260.3580701710861
$ws.Range("S21").Value = 0.002710670808468793
$ws.Range("T21").Value = 0.002710670808468793
$ws.Range("G22").Value = 16.16161433333333
$ws.Range("H22").Value = 48.484843
$ws.Range("I22").Value = 0.01048312291409786
$ws.Range("J22").Value = 0.01048312291409786
$ws.Range("M22").Value = 1.780241
$ws.Range("N22").Value = 5.340723000000001
$ws.Range("O22").Value = 0.2571704734300857
$ws.Range("P22").Value = 0.2571704734300857
$ws.Range("Q22").Value = 28.77156846238767
$ws.Range("R22").Value = 258.944116161489
$ws.Range("S22").Value = 0.002695949682844328
$ws.Range("T22").Value = 0.002695949682844327
$ws.Range("G23").Value = 16.16161433333333
$ws.Range("H23").Value = 48.484843
$ws.Range("I23").Value = 0.01048312291409786
$ws.Range("J23").Value = 0.01048312291409786
$ws.Range("O23").Value = 0.4607709215973151
$ws.Range("P23").Value = 0.4607709215973152
$ws.Range("Q23").Value = 51.549860835088
$ws.Range("R23").Value = 463.948747515792
$ws.Range("S23").Value = 0.004830318206346805
$ws.Range("T23").Value = 0.004830318206346804
$ws.Range("G24").Value = 16.16161433333333
$ws.Range("H24").Value = 48.484843
$ws.Range("I24").Value = 0.01048312291409786
$ws.Range("J24").Value = 0.01048312291409786
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.4746316666666666
$ws.Range("N24").Value = 1.423895
$ws.Range("O24").Value = 0.06856445302718973
$ws.Range("P24").Value = 0.06856445302718973
$ws.Range("Q24").Value = 7.670813947053888
$ws.Range("R24").Value = 69.03732552348499
$ws.Range("S24").Value = 0.0007187695886219194
$ws.Range("T24").Value = 0.0007187695886219193
$ws.Range("G25").Value = 16.16161433333333
$ws.Range("H25").Value = 48.484843
$ws.Range("I25").Value = 0.01048312291409786
$ws.Range("J25").Value = 0.01048312291409786
$ws.Range("M25").Value = 0.9303213333333334
$ws.Range("N25").Value = 2.790964
$ws.Range("O25").Value = 0.1343925781596098
$ws.Range("P25").Value = 0.1343925781596098
$ws.Range("Q25").Value = 15.03549459540578
$ws.Range("R25").Value = 135.319451358652
$ws.Range("S25").Value = 0.001408853915589694
$ws.Range("T25").Value = 0.001408853915589694
$ws.Range("G26").Value = 16.16161433333333
$ws.Range("H26").Value = 48.484843
$ws.Range("I26").Value = 0.01048312291409786
$ws.Range("J26").Value = 0.01048312291409786
$ws.Range("M26").Value = 0.547574
$ws.Range("N26").Value = 1.642722
$ws.Range("O26").Value = 0.07910157378579964
$ws.Range("P26").Value = 0.07910157378579964
$ws.Range("Q26").Value = 8.849679806960667
$ws.Range("R26").Value = 79.647118262646
$ws.Range("S26").Value = 0.0008292315206951192
$ws.Range("T26").Value = 0.000829231520695119
